$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old columns C and D entirely (headers + data)
$ws.Range("C1:D7").Clear()

# Update header row
$ws.Range("B1").Value = "New_solver.py"

# Update data rows with new test names and timings
$data = @(
    @("small_test_1.txt", 3.724148035049438),
    @("small_test_2.txt", 1.91952109336853),
    @("small_test_3.txt", 3.322933197021484),
    @("medium_test_1.txt", 119.0632381439209),
    @("medium_test_2.txt", 67.22225594520569),
    @("medium_test_3.txt", 9.359437942504883),
    @("large_test_2.txt", 57.22832894325256),
    @("large_test_3.txt", 54.76760601997375)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
